# ---------------------------------------------------------------------------
# "it now reads/writes to same file"
#
# The TA (trailing-average) summary that used to be read off Sheet1's time
# series and written to Sheet2 is now written to Sheet3 instead (Sheet2 is
# left blank), and the results reflect the new window grabbed from Sheet1
# (5 rows -> 3 rows, recomputed D/E columns). A new, still-blank tab
# (Sheet4) is appended and becomes the active/selected sheet - the "write
# the TA to a different tab" step next time round.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Sheet2 no longer gets written to -> clear its old 5-row summary back to
# a blank sheet.
$ws2.Range("A1:E5").ClearContents()

# Sheet3 now receives the (re-windowed) 3-row summary.
$ws3.Cells.Item(1, 1).Value = 432.9604174527465
$ws3.Cells.Item(1, 2).Value = 157.6565028933628
$ws3.Cells.Item(1, 3).Value = 100.5364509608477
$ws3.Cells.Item(1, 4).Value = 19069.94516042675
$ws3.Cells.Item(1, 5).Value = 31490.55518780852

$ws3.Cells.Item(2, 1).Value = 385.142419893523
$ws3.Cells.Item(2, 2).Value = 129.2841508036708
$ws3.Cells.Item(2, 3).Value = 108.1996807442682
$ws3.Cells.Item(2, 4).Value = 15891.04308403571
$ws3.Cells.Item(2, 5).Value = 26805.16640513378

$ws3.Cells.Item(3, 1).Value = 409.0514186731344
$ws3.Cells.Item(3, 2).Value = 143.7454045971708
$ws3.Cells.Item(3, 3).Value = 104.9459852064982
$ws3.Cells.Item(3, 4).Value = 17480.49412223122
$ws3.Cells.Item(3, 5).Value = 29147.8607964711

# Append the new destination tab (Sheet4) after Sheet3, matching outline/
# page-setup of the existing sheets, then make it the active tab/selection
# (activeTab 0 -> 3, selection E28 on Sheet4).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet4"
$ws4.Outline.SummaryRow = 1
$ws4.Outline.SummaryColumn = 1

$ws4.Activate()
$ws4.Range("E28").Select()
